$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Formula changes: move from an incremental (row-to-row delta) stress/strain
# formulation to a total stress/strain formulation.
# ---------------------------------------------------------------------------

# AI column: -AG/AH (deviator stress ratio using AG/AH deltas) -> -AD/AE (total)
$ws.Range("AI3").Formula = "=-AD3/AE3"
$ws.Range("AI4:AI18").Formula = "=-AD4/AE4"

# AP column: secant-slope-between-rows ratio -> direct total ratio via AJ^0.25
$ws.Range("AP3").Formula = "=AJ3^0.25"
$ws.Range("AP4:AP18").Formula = "=AJ4^0.25"

# AS column: cumulative running sum (AS2+AR3, AS3+AR4, ...) -> total-form
# expression evaluated one row ahead (same pattern AR/AT already used, but
# without the "current minus previous" delta and referencing the next row).
$ws.Range("AS2").Formula = "=(1+2*AM3)*AK3*(1-AP3/3)/(2*AN3*AO3*AM3)"
$ws.Range("AS3:AS18").Formula = "=(1+2*AM4)*AK4*(1-AP4/3)/(2*AN4*AO4*AM4)"

# AU column: cumulative running sum (AU2+AT3, AU3+AT4, ...) -> total-form
# expression (same structure as AT, but using AK instead of AK-previous AK).
$ws.Range("AU3").Formula = "=(1-AM3)*(AK3)*(1-AP3/3)/(3*AN3*AO3*AM3)"
$ws.Range("AU4:AU18").Formula = "=(1-AM4)*(AK4)*(1-AP4/3)/(3*AN4*AO4*AM4)"

# ---------------------------------------------------------------------------
# Formatting: the newly-formula-driven cells AI4:AI18, AP4:AP18, AS2,
# AS4:AS18 and AU4:AU18 pick up the yellow "input/highlight" fill that the
# rest of the AI/AP/AS/AU/AO columns already use (row 3 already carries that
# fill via its row-level style).
# ---------------------------------------------------------------------------
$highlight = $ws.Range("AO4").Interior.ColorIndex
$ws.Range("AS2").Interior.ColorIndex = $highlight
$ws.Range("AI4:AI18").Interior.ColorIndex = $highlight
$ws.Range("AP4:AP18").Interior.ColorIndex = $highlight
$ws.Range("AS4:AS18").Interior.ColorIndex = $highlight
$ws.Range("AU4:AU18").Interior.ColorIndex = $highlight

# ---------------------------------------------------------------------------
# View state: zoom out, scroll the window over one column and move the
# live "entire column" selection from AM to AO.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 55
$win.ScrollColumn = 28
$ws.Range("AO1:AO1048576").Select()
